$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append 10 new rows (147-156) of test data for the reg_center_device
# master table. Columns:
#   A regcntr_id   (constant 10001)
#   B device_id    (incrementing MAC/device id)
#   C lang_code    (eng)
#   D is_active    (TRUE)
#   E cr_by        (superadmin)
#   F cr_dtimes    (now())

$startRow = 147
$startDevice = 3000166
$count = 10

for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 10001
    $ws.Cells.Item($r, 2).Value = $startDevice + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

$ws.Range("C152").Select()
